$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# WAT-413 (row 50): clarify the description now that a sibling scenario (sorting by
# totalNumberOfPublications) exists - this one specifically covers sorting using years.
$ws.Range("B50").Value = "Verify that user should be able to filter values for a set of authors  provided with all the mandatory inputs along with sorting using years given an order"

# New WAT-415 test case (row 51): same author-search scenario, but sorted/validated by
# totalNumberOfPublications instead of publication year.
$ws.Range("A51").Value = "WAT-415"
$ws.Range("B51").Value = "Verify that user should be able to filter values for a set of authors provided with all the mandatory inputs along with sorting using totalNumberOfPublications given an order"
$ws.Range("C51").Value = "WOSAUTHORRECOMMEND"
$ws.Range("D51").Value = "/author/search"
$ws.Range("E51").Value = "GET"
$ws.Range("G51").Value = "?name=wang&affiliation=china&filter=name&category=physics&filter=affiliation&filter=catagory&sort=totalNumberOfPublications&order=asc&limit=10"
$ws.Range("J51").Value = "status=200||hits.primaryName=wang||filters.category=physics||hits.affiliation=china"
$ws.Range("K51").Value = "hits[0].totalNumberOfPublications"

# Carry over row 50's look (borders/alignment/wrap/row height) onto the newly added row 51.
$ws.Range("A50:L50").Copy()
$ws.Range("A51:L51").PasteSpecial(-4122)
$ws.Range("A51:L51").RowHeight = 45

$ws.Range("K51").Select()
